$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Two test-case names picked up a "(-)" suffix.
$ws.Range("B4").Value = "Submitting Empty Fields(-)"
$ws.Range("B20").Value = "No Results Found(-)"

# View tweaks: zoom to 85%, drop the old scrolled/top-left position, and
# move the live selection onto B20.
$ws.Activate()
$excel.ActiveWindow.Zoom = 85
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B20").Select()
